# Add a new column J ("CatCore Team") to the TFVC permissions report,
# mirroring the existing "CatCore Build Service (DefaultCollection)" column (I),
# and widen a few columns to accommodate the longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns to match new layout (values tuned so the engine's internal
# char-width snapping lands as close as possible to the target stored widths).
$ws.Columns.Item(3).ColumnWidth = 18.833333333333336
$ws.Columns.Item(5).ColumnWidth = 23.833333333333336
$ws.Columns.Item(8).ColumnWidth = 42.83333333333333
$ws.Columns.Item(10).ColumnWidth = 42.83333333333333

# Update "Allow" -> "Allow (inherited)" across the permission grid (columns B-I, rows 2-14).
for ($r = 2; $r -le 14; $r++) {
    for ($c = 2; $c -le 9; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cellText = $cell.Text
        if ($cellText -eq "Allow") {
            $cell.Value = "Allow (inherited)"
        }
    }
}

# New header for column J, duplicating column I's header text.
$headerText = $ws.Range("I1").Text
$ws.Range("J1").Value = $headerText
$ws.Range("J1").Font.Bold = $true

# Fill column J (rows 2-14) with the same value as column I for that row
# (after the Allow -> Allow (inherited) update above).
for ($r = 2; $r -le 14; $r++) {
    $iText = $ws.Cells.Item($r, 9).Text
    $ws.Cells.Item($r, 10).Value = $iText
}
